$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$queryCases = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE f.file_format IN ['bai'] 
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE f.file_format IN ['bai'] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

$queryFiles = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE f.file_format IN ['bai'] 
 WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# Insert a new column A for "TabName" labels: shift existing columns A:D to B:E first.
$ws.Columns.Item(1).Insert()

# Row 1 header
$ws.Cells.Item(1,1).Value = "TabName"

# Row 2 (Cases tab)
$ws.Cells.Item(2,1).Value = "CasesTab"
$ws.Cells.Item(2,2).Value = $queryCases
$ws.Cells.Item(2,3).Value = $statQuery

# Row 3 (Files tab) - new row
$ws.Cells.Item(3,1).Value = "FilesTab"
$ws.Cells.Item(3,2).Value = $queryFiles
$ws.Cells.Item(3,3).Value = $statQuery
$ws.Cells.Item(3,4).Value = "TC01_Trials_Filter_AssocFileFormat-Bai_Neo4jData.xlsx"
$ws.Cells.Item(3,5).Value = "TC01_Trials_Filter_AssocFileFormat-Bai_WebData.xlsx"

# Apply the wrap-text style (style index 1, "Normal 2") to the query cells, matching B2/C2/B3/C3
$ws.Cells.Item(2,2).WrapText = $true
$ws.Cells.Item(2,3).WrapText = $true
$ws.Cells.Item(3,2).WrapText = $true
$ws.Cells.Item(3,3).WrapText = $true

# Column A is brand new - give it a narrow width close to the target 8.81640625
# (columns B:E already inherited their correct widths from the original A:D via the
# Insert() shift above, so we deliberately leave them untouched to avoid re-rounding them).
$ws.Columns.Item(1).ColumnWidth = 8

# Row heights
$ws.Rows.Item(2).RowHeight = 188.5
$ws.Rows.Item(3).RowHeight = 409.5

# Sheet view settings
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("B2").Select()
